# Update the two-digit division worksheet numbers.
# The worksheet is a single table; every 4th row (1, 5, 9, 13, 17) holds
# five division problems (columns 1-5). Replace each cell's text with its
# new value, addressing cells by (row, column) so duplicate old values
# (e.g. "78÷4=" appears twice) are each mapped to their own new value.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="37÷4="},
    @{Row=1;  Col=2; Text="73÷9="},
    @{Row=1;  Col=3; Text="14÷9="},
    @{Row=1;  Col=4; Text="43÷4="},
    @{Row=1;  Col=5; Text="66÷8="},

    @{Row=5;  Col=1; Text="19÷2="},
    @{Row=5;  Col=2; Text="29÷4="},
    @{Row=5;  Col=3; Text="49÷5="},
    @{Row=5;  Col=4; Text="87÷9="},
    @{Row=5;  Col=5; Text="34÷2="},

    @{Row=9;  Col=1; Text="97÷9="},
    @{Row=9;  Col=2; Text="35÷8="},
    @{Row=9;  Col=3; Text="52÷5="},
    @{Row=9;  Col=4; Text="62÷9="},
    @{Row=9;  Col=5; Text="41÷8="},

    @{Row=13; Col=1; Text="84÷5="},
    @{Row=13; Col=2; Text="17÷3="},
    @{Row=13; Col=3; Text="29÷6="},
    @{Row=13; Col=4; Text="84÷8="},
    @{Row=13; Col=5; Text="34÷4="},

    @{Row=17; Col=1; Text="62÷9="},
    @{Row=17; Col=2; Text="79÷3="},
    @{Row=17; Col=3; Text="37÷2="},
    @{Row=17; Col=4; Text="43÷9="},
    @{Row=17; Col=5; Text="51÷9="}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
